$wb = $excel.ActiveWorkbook

# --- zh-cn sheet (row 7: 7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb) ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.0e67ba9e1a0480b0aba07d041b05efa2efb18c9a.zh-cn.xlf"
$wsZh.Range("K7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3d0ab96f89eb8191ef2ba4ca5e5cc48b7590ec39/e2e/7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7e190263c4dce59ad5553f87ae787bf5db13e98/e2e/7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.md."
$wsZh.Range("P7").Value = "2016-08-30 04:54:36"

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/3d0ab96f89eb8191ef2ba4ca5e5cc48b7590ec39/e2e/7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.md", "", "", "7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.md")

# --- de-de sheet (row 7: 7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb) ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.0e67ba9e1a0480b0aba07d041b05efa2efb18c9a.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-30 04:55:09"
$wsDe.Range("P7").Value = "2016-08-30 04:54:36"

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3d0ab96f89eb8191ef2ba4ca5e5cc48b7590ec39/e2e/7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.md", "", "", "7dca80d8-a6bf-4756-a4c9-b6cf0ff8aadb.md")
